$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1923076923076923
$ws.Range("C2").Value = 0.5629370629370629
$ws.Range("J2").Value = 0.01748251748251748
$ws.Range("O2").Value = 0.003496503496503497
$ws.Range("P2").Value = 0.1293706293706294
$ws.Range("S2").Value = 0.0944055944055944
$ws.Range("B3").Value = 0.01197604790419162
$ws.Range("C3").Value = 0.02994011976047904
$ws.Range("J3").Value = 0.005988023952095809
$ws.Range("P3").Value = 0.7724550898203593
$ws.Range("S3").Value = 0.1796407185628743
$ws.Range("B6").Value = 0.04868913857677903
$ws.Range("D6").Value = 0.00749063670411985
$ws.Range("F6").Value = 0.05617977528089887
$ws.Range("J6").Value = 0.2808988764044944
$ws.Range("O6").Value = 0.003745318352059925
$ws.Range("Q6").Value = 0.1161048689138577
$ws.Range("R6").Value = 0.1161048689138577
$ws.Range("S6").Value = 0.3707865168539326
$ws.Range("B7").Value = 0.09482758620689655
$ws.Range("D7").Value = 0.02586206896551724
$ws.Range("F7").Value = 0.04310344827586207
$ws.Range("J7").Value = 0.09913793103448276
$ws.Range("O7").Value = 0.02586206896551724
$ws.Range("Q7").Value = 0.1982758620689655
$ws.Range("R7").Value = 0.05603448275862069
$ws.Range("S7").Value = 0.456896551724138
$ws.Range("B8").Value = 0.0797872340425532
$ws.Range("D8").Value = 0.008865248226950355
$ws.Range("F8").Value = 0.0549645390070922
$ws.Range("J8").Value = 0.1400709219858156
$ws.Range("O8").Value = 0.0124113475177305
$ws.Range("Q8").Value = 0.173758865248227
$ws.Range("R8").Value = 0.09397163120567376
$ws.Range("S8").Value = 0.4361702127659575
$ws.Range("B9").Value = 0.1038251366120219
$ws.Range("F9").Value = 0.09836065573770492
$ws.Range("J9").Value = 0.1311475409836066
$ws.Range("O9").Value = 0.01092896174863388
$ws.Range("Q9").Value = 0.1420765027322404
$ws.Range("R9").Value = 0.08743169398907104
$ws.Range("S9").Value = 0.4262295081967213
$ws.Range("B10").Value = 0.08605135322692574
$ws.Range("D10").Value = 0.02498265093684941
$ws.Range("E10").Value = 0.001387925052047189
$ws.Range("F10").Value = 0.06731436502428868
$ws.Range("J10").Value = 0.1304649548924358
$ws.Range("O10").Value = 0.007633587786259542
$ws.Range("Q10").Value = 0.2116585704371964
$ws.Range("R10").Value = 0.08188757807078417
$ws.Range("S10").Value = 0.3886190145732131
$ws.Range("G11").Value = 0.1195652173913044
$ws.Range("J11").Value = 0.1032608695652174
$ws.Range("K11").Value = 0.1766304347826087
$ws.Range("L11").Value = 0.5733695652173914
$ws.Range("S11").Value = 0.02717391304347826
$ws.Range("G12").Value = 0.7568807339449541
$ws.Range("J12").Value = 0.1605504587155963
$ws.Range("K12").Value = 0.004587155963302753
$ws.Range("L12").Value = 0.01376146788990826
$ws.Range("S12").Value = 0.06422018348623854
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.2692307692307692
$ws.Range("S13").Value = 0.01923076923076923
$ws.Range("F15").Value = 0.02643171806167401
$ws.Range("H15").Value = 0.2114537444933921
$ws.Range("I15").Value = 0.06607929515418502
$ws.Range("J15").Value = 0.3612334801762114
$ws.Range("K15").Value = 0.06167400881057269
$ws.Range("M15").Value = 0.02202643171806168
$ws.Range("O15").Value = 0.06167400881057269
$ws.Range("S15").Value = 0.1894273127753304
$ws.Range("F16").Value = 0.02083333333333333
$ws.Range("H16").Value = 0.1822916666666667
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.3854166666666667
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.03125
$ws.Range("N16").Value = 0.005208333333333333
$ws.Range("O16").Value = 0.07291666666666667
$ws.Range("S16").Value = 0.15625
$ws.Range("F17").Value = 0.02994011976047904
$ws.Range("H17").Value = 0.1976047904191617
$ws.Range("I17").Value = 0.07385229540918163
$ws.Range("J17").Value = 0.3992015968063872
$ws.Range("K17").Value = 0.1197604790419162
$ws.Range("M17").Value = 0.01796407185628742
$ws.Range("O17").Value = 0.0439121756487026
$ws.Range("S17").Value = 0.1177644710578842
$ws.Range("F18").Value = 0.02597402597402598
$ws.Range("H18").Value = 0.1904761904761905
$ws.Range("I18").Value = 0.08658008658008658
$ws.Range("J18").Value = 0.3593073593073593
$ws.Range("K18").Value = 0.1341991341991342
$ws.Range("M18").Value = 0.01298701298701299
$ws.Range("O18").Value = 0.06060606060606061
$ws.Range("S18").Value = 0.1298701298701299
$ws.Range("F19").Value = 0.01954177897574124
$ws.Range("H19").Value = 0.2257412398921833
$ws.Range("I19").Value = 0.06873315363881402
$ws.Range("J19").Value = 0.3605121293800539
$ws.Range("K19").Value = 0.1206199460916442
$ws.Range("M19").Value = 0.01954177897574124
$ws.Range("N19").Value = 0.0006738544474393531
$ws.Range("O19").Value = 0.07749326145552561
$ws.Range("S19").Value = 0.1071428571428571
